$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 29800
$ws.Range("J21").Value = 0
$ws.Range("L21").Value = 0
$ws.Range("N21").ClearContents()
$ws.Range("H23").Value = 29800
$ws.Range("J23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("N23").ClearContents()
$ws.Range("H41").Value = 358.57144
$ws.Range("I41").Value = 164.2
$ws.Range("J41").Value = 844.5
$ws.Range("K41").Value = 164.2
$ws.Range("L41").Value = 844.5
$ws.Range("M41").Value = 275.8
$ws.Range("N41").Value = -1724.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 4593.3
$ws.Range("I28").Value = 4593.3
$ws.Range("K28").Value = 4593.3
$ws.Range("M28").Value = -4401.3
$ws.Range("H32").Value = 3740.42
$ws.Range("I32").Value = 2854.975
$ws.Range("K32").Value = 2854.975
$ws.Range("M32").Value = -2567.975
$ws.Range("H36").Value = 9500
$ws.Range("I36").Value = 9500
$ws.Range("K36").Value = 9500
$ws.Range("M36").Value = -9154
$ws.Range("H61").Value = 2170.5918
$ws.Range("I61").Value = 1696.5938
$ws.Range("K61").Value = 1696.5938
$ws.Range("M61").Value = -1484.5938
$ws.Range("H74").Value = 194095.69
$ws.Range("I74").Value = 371655
$ws.Range("K74").Value = 371655
$ws.Range("M74").Value = -370781
$ws.Range("H77").Value = 194095.69
$ws.Range("I77").Value = 371655
$ws.Range("K77").Value = 1858275
$ws.Range("M77").Value = -1853907
$ws.Range("H99").Value = 4593.3
$ws.Range("I99").Value = 4593.3
$ws.Range("K99").Value = 4593.3
$ws.Range("M99").Value = -1598.3
$ws.Range("H110").Value = 2931.3914
$ws.Range("I110").Value = 1590.4375
$ws.Range("K110").Value = 1590.4375
$ws.Range("M110").Value = 454.5625
$ws.Range("H132").Value = 1684.2667
$ws.Range("I132").Value = 894.3714
$ws.Range("J132").Value = 4448.9
$ws.Range("K132").Value = 2683.1142
$ws.Range("L132").Value = 13346.7
$ws.Range("M132").Value = -153.1142
$ws.Range("N132").Value = -18406.7
$ws.Range("H135").Value = 112954.2
$ws.Range("J135").Value = 112954.2
$ws.Range("L135").Value = 112954.2
$ws.Range("N135").Value = -123094.2
$ws.Range("H136").Value = 2170.5918
$ws.Range("I136").Value = 1696.5938
$ws.Range("K136").Value = 5089.7814
$ws.Range("M136").Value = -2539.7814

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 7776.25
$ws.Range("I26").Value = 7776.25
$ws.Range("K26").Value = 7776.25
$ws.Range("M26").Value = -7484.25
$ws.Range("H82").Value = 41462
$ws.Range("J82").Value = 81420
$ws.Range("L82").Value = 81420
$ws.Range("N82").Value = -82186
$ws.Range("H85").Value = 41462
$ws.Range("J85").Value = 81420
$ws.Range("L85").Value = 81420
$ws.Range("N85").Value = -84072
$ws.Range("H102").Value = 0
$ws.Range("I102").Value = 0
$ws.Range("K102").Value = 0
$ws.Range("M102").ClearContents()
$ws.Range("H105").Value = 8669194
$ws.Range("I105").Value = 502497.34
$ws.Range("K105").Value = 502497.34
$ws.Range("M105").Value = -500750.34
$ws.Range("H107").Value = 1604.8572
$ws.Range("I107").Value = 1430.2632
$ws.Range("J107").Value = 1973.4445
$ws.Range("K107").Value = 1430.2632
$ws.Range("L107").Value = 1973.4445
$ws.Range("M107").Value = 489.7367999999999
$ws.Range("N107").Value = -5813.4445
$ws.Range("H134").Value = 3239.0676
$ws.Range("I134").Value = 3013.2068
$ws.Range("J134").Value = 4057.8125
$ws.Range("K134").Value = 9039.6204
$ws.Range("L134").Value = 12173.4375
$ws.Range("M134").Value = -6504.6204
$ws.Range("N134").Value = -17243.4375

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 1086.4445
$ws.Range("I22").Value = 1870
$ws.Range("J22").Value = 694.6667
$ws.Range("K22").Value = 1870
$ws.Range("L22").Value = 694.6667
$ws.Range("M22").Value = -1520
$ws.Range("N22").Value = -1394.6667
$ws.Range("H122").Value = 3338.36
$ws.Range("I122").Value = 2314.6667
$ws.Range("K122").Value = 6944.000100000001
$ws.Range("M122").Value = -4494.000100000001
$ws.Range("H132").Value = 4137.0713
$ws.Range("I132").Value = 4347.7144
$ws.Range("K132").Value = 13043.1432
$ws.Range("M132").Value = -10513.1432

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H100").Value = 5408.6665
$ws.Range("J100").Value = 5408.6665
$ws.Range("L100").Value = 16225.9995
$ws.Range("N100").Value = -17847.9995
$ws.Range("H131").Value = 2187.9707
$ws.Range("I131").Value = 2220.8462
$ws.Range("J131").Value = 2167.6191
$ws.Range("K131").Value = 6662.5386
$ws.Range("L131").Value = 6502.8573
$ws.Range("M131").Value = -1622.5386
$ws.Range("N131").Value = -16582.8573

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2869.5386
$ws.Range("I132").Value = 2347.95
$ws.Range("J132").Value = 4608.1665
$ws.Range("K132").Value = 7043.849999999999
$ws.Range("L132").Value = 13824.4995
$ws.Range("M132").Value = -4513.849999999999
$ws.Range("N132").Value = -18884.4995

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H99").Value = 24984
$ws.Range("I99").Value = 24984
$ws.Range("K99").Value = 24984
$ws.Range("M99").Value = -21989
$ws.Range("H122").Value = 2930.9
$ws.Range("I122").Value = 2929.1
$ws.Range("J122").Value = 2932.7
$ws.Range("K122").Value = 8787.299999999999
$ws.Range("L122").Value = 8798.099999999999
$ws.Range("M122").Value = -6337.299999999999
$ws.Range("N122").Value = -13698.1
$ws.Range("H136").Value = 5503.077
$ws.Range("I136").Value = 5282.3335
$ws.Range("K136").Value = 15847.0005
$ws.Range("M136").Value = -13297.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H52").Value = 22997.2
$ws.Range("I52").Value = 18332.334
$ws.Range("K52").Value = 18332.334
$ws.Range("M52").Value = -18106.334
$ws.Range("H58").Value = 18553
$ws.Range("J58").Value = 0
$ws.Range("L58").Value = 0
$ws.Range("N58").ClearContents()
$ws.Range("H122").Value = 11364944
$ws.Range("I122").Value = 1283.6
$ws.Range("J122").Value = 35715644
$ws.Range("K122").Value = 3850.8
$ws.Range("L122").Value = 107146932
$ws.Range("M122").Value = -1400.8
$ws.Range("N122").Value = -107151832
$ws.Range("H126").Value = 4578.2
$ws.Range("I126").Value = 1191
$ws.Range("J126").Value = 5425
$ws.Range("K126").Value = 3573
$ws.Range("L126").Value = 16275
$ws.Range("M126").Value = -1103
$ws.Range("N126").Value = -21215
$ws.Range("H136").Value = 83338550
$ws.Range("I136").Value = 100001144
$ws.Range("J136").Value = 25625
$ws.Range("K136").Value = 300003432
$ws.Range("L136").Value = 76875
$ws.Range("M136").Value = -300000882
$ws.Range("N136").Value = -81975
